$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# GanttChart correction: the "Matches Management" task start date (row 7)
# moves from 3/5/2018 (serial 43164) to 3/7/2018 (serial 43166).
$ws.Range("B7").Value = 43166

# The sheet's last saved selection moves from E8 to B7.
$ws.Activate()
$ws.Range("B7").Select()
